# Generate Report for Handoff
# Replace the old GUID-based file identifiers / hashes and timestamps with the
# new ones produced by this handoff run.

$wb = $excel.ActiveWorkbook

$oldGuid = "342141b0-ae3a-4a0d-9acd-9d4a64b23197"
$newGuid = "3fec1a29-f8dc-499b-adbd-64b948868ea9"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-21 17:03:52"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.f5863a18e1d03ab84f804234215253f1aa5fb653.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-21 17:03:49"
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.f5863a18e1d03ab84f804234215253f1aa5fb653.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-21 17:03:52"
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
